$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.583.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.94%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.626.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.72"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +25.33%  "

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "224.85"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "644.48"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.424"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -4.05%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.622.46"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.73"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +14.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.219"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000296"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.53"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.304.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.49%  "

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +31.48%  "

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.375.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.10"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.92"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.622.27"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.301"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +37.98%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "136.21"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +13.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "533.24"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.31"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.58%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.70%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000204"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.36"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.796.26"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.63"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.18"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.31%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.647"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.89"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "33.69"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.61%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.17%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +21.81%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.38%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "599.16"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.33%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.27"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.511"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.68%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.23"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "237.93"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +13.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.36"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.64%  "
